$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 7.132470666666666
$ws.Range("H2").Value = 21.397412
$ws.Range("I2").Value = 0.1078130252899183
$ws.Range("J2").Value = 0.1078130252899183
$ws.Range("M2").Value = 4.47329
$ws.Range("N2").Value = 13.41987
$ws.Range("O2").Value = 0.3468876470949054
$ws.Range("P2").Value = 0.3468876470949054
$ws.Range("Q2").Value = 31.90560970849333
$ws.Range("R2").Value = 287.15048737644
$ws.Range("S2").Value = 0.0373990066690033
$ws.Range("T2").Value = 0.0373990066690033

# Row 3
$ws.Range("G3").Value = 7.132470666666666
$ws.Range("H3").Value = 21.397412
$ws.Range("I3").Value = 0.1078130252899183
$ws.Range("J3").Value = 0.1078130252899183
$ws.Range("O3").Value = 0.3372845821706862
$ws.Range("P3").Value = 0.3372845821706862
$ws.Range("Q3").Value = 31.02235069352577
$ws.Range("R3").Value = 279.201156241732
$ws.Range("S3").Value = 0.03636367118746774
$ws.Range("T3").Value = 0.03636367118746774

# Row 4
$ws.Range("G4").Value = 7.132470666666666
$ws.Range("H4").Value = 21.397412
$ws.Range("I4").Value = 0.1078130252899183
$ws.Range("J4").Value = 0.1078130252899183
$ws.Range("M4").Value = 4.072757333333333
$ws.Range("N4").Value = 12.218272
$ws.Range("O4").Value = 0.3158277707344083
$ws.Range("P4").Value = 0.3158277707344083
$ws.Range("Q4").Value = 29.04882221245155
$ws.Range("R4").Value = 261.4393999120639
$ws.Range("S4").Value = 0.03405034743344729
$ws.Range("T4").Value = 0.03405034743344729

# Row 5
$ws.Range("I5").Value = 0.2490596131114117
$ws.Range("J5").Value = 0.2490596131114118
$ws.Range("M5").Value = 4.47329
$ws.Range("N5").Value = 13.41987
$ws.Range("O5").Value = 0.3468876470949054
$ws.Range("P5").Value = 0.3468876470949054
$ws.Range("Q5").Value = 73.70536898220334
$ws.Range("R5").Value = 663.3483208398301
$ws.Range("S5").Value = 0.08639570317858505
$ws.Range("T5").Value = 0.08639570317858507

# Row 6
$ws.Range("I6").Value = 0.2490596131114117
$ws.Range("J6").Value = 0.2490596131114118
$ws.Range("O6").Value = 0.3372845821706862
$ws.Range("P6").Value = 0.3372845821706862
$ws.Range("S6").Value = 0.08400396754387528
$ws.Range("T6").Value = 0.08400396754387528

# Row 7
$ws.Range("I7").Value = 0.2490596131114117
$ws.Range("J7").Value = 0.2490596131114118
$ws.Range("M7").Value = 4.072757333333333
$ws.Range("N7").Value = 12.218272
$ws.Range("O7").Value = 0.3158277707344083
$ws.Range("P7").Value = 0.3158277707344083
$ws.Range("Q7").Value = 67.10588448956089
$ws.Range("R7").Value = 603.9529604060481
$ws.Range("S7").Value = 0.07865994238895137
$ws.Range("T7").Value = 0.07865994238895137

# Row 8
$ws.Range("G8").Value = 42.546687
$ws.Range("H8").Value = 127.640061
$ws.Range("I8").Value = 0.6431273615986699
$ws.Range("J8").Value = 0.6431273615986699
$ws.Range("M8").Value = 4.47329
$ws.Range("N8").Value = 13.41987
$ws.Range("O8").Value = 0.3468876470949054
$ws.Range("P8").Value = 0.3468876470949054
$ws.Range("Q8").Value = 190.32366949023
$ws.Range("R8").Value = 1712.91302541207
$ws.Range("S8").Value = 0.223092937247317
$ws.Range("T8").Value = 0.223092937247317

# Row 9
$ws.Range("G9").Value = 42.546687
$ws.Range("H9").Value = 127.640061
$ws.Range("I9").Value = 0.6431273615986699
$ws.Range("J9").Value = 0.6431273615986699
$ws.Range("O9").Value = 0.3372845821706862
$ws.Range("P9").Value = 0.3372845821706862
$ws.Range("Q9").Value = 185.054843776669
$ws.Range("R9").Value = 1665.493593990021
$ws.Range("S9").Value = 0.2169169434393432
$ws.Range("T9").Value = 0.2169169434393432

# Row 10
$ws.Range("G10").Value = 42.546687
$ws.Range("H10").Value = 127.640061
$ws.Range("I10").Value = 0.6431273615986699
$ws.Range("J10").Value = 0.6431273615986699
$ws.Range("M10").Value = 4.072757333333333
$ws.Range("N10").Value = 12.218272
$ws.Range("O10").Value = 0.3158277707344083
$ws.Range("P10").Value = 0.3158277707344083
$ws.Range("Q10").Value = 173.282331488288
$ws.Range("R10").Value = 1559.540983394592
$ws.Range("S10").Value = 0.2031174809120096
$ws.Range("T10").Value = 0.2031174809120096
